$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.959.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.18%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.570.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.62%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.10%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '568.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.23%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.562.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.74%  '

$ws.Range('E8').Value = '  -1.38%  '

$ws.Range('E9').Value = '  +0.24%  '

$ws.Range('E10').Value = '  -0.74%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '61.60'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.72%  '

$ws.Range('E12').Value = '  -2.71%  '

$ws.Range('E13').Value = '  +3.38%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.19'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.79%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.132.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.68%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.554.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.82%  '

$ws.Range('E17').Value = '  +0.40%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '67.666.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.51%  '

$ws.Range('E21').Value = '  -0.77%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '400.60'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.43%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.63%  '

$ws.Range('E24').Value = '  -1.72%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.04'
$ws.Range('D25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.47%  '

$ws.Range('E27').Value = '  -2.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.41%  '

$ws.Range('E30').Value = '  -1.57%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.83%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '669.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.02'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.72%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.18'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.02%  '

$ws.Range('E35').Value = '  -2.43%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '40.75'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.85%  '

$ws.Range('E37').Value = '  +0.27%  '

$ws.Range('E38').Value = '  -0.10%  '

$ws.Range('E39').Value = '  +9.56%  '

$ws.Range('E40').Value = '  -2.78%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.158.57'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.14%  '

$ws.Range('E42').Value = '  -0.90%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.997'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.20%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.64'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.91%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.11%  '

$ws.Range('E46').Value = '  +11.32%  '

$ws.Range('E47').Value = '  -1.78%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.130'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.01%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.00%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '137.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.92%  '
